$d = $word.ActiveDocument

# 1) "Achelin Felix - Front-end Developer - " (removes proofErr spell-check
#    split around "Achelin" by merging the two runs back into one)
$d.Content.Find.Execute("Achelin Felix - Front-end Developer - ", $false, $false, $false, $false, $false, $true, 1, $false, "Achelin Felix - Front-end Developer - ", 2)

# 2) "Samuel Yambo -Front/Back-end Developer- " (removes proofErr around "Yambo")
$d.Content.Find.Execute("Samuel Yambo -Front/Back-end Developer- ", $false, $false, $false, $false, $false, $true, 1, $false, "Samuel Yambo -Front/Back-end Developer- ", 2)

# 3) "Frecks Bertrand II – Front-end Developer - " (removes proofErr around "Frecks")
$d.Content.Find.Execute("Frecks Bertrand II " + [char]0x2013 + " Front-end Developer - ", $false, $false, $false, $false, $false, $true, 1, $false, "Frecks Bertrand II " + [char]0x2013 + " Front-end Developer - ", 2)

# 4) Remove " (may delete)" from the Homepage bullet
$d.Content.Find.Execute("personal account page (may delete) ", $false, $false, $false, $false, $false, $true, 1, $false, "personal account page ", 2)

# 5) "...footer identifies you as "testUser" ..." (removes proofErr around "testUser")
$d.Content.Find.Execute("The forum homepage will display all topics and a preview of posts within those topics. Please ensure that the footer identifies you as " + [char]0x201C + "testUser" + [char]0x201D + " and that the data on the forum homepage is useful and accurate.", $false, $false, $false, $false, $false, $true, 1, $false, "The forum homepage will display all topics and a preview of posts within those topics. Please ensure that the footer identifies you as " + [char]0x201C + "testUser" + [char]0x201D + " and that the data on the forum homepage is useful and accurate.", 2)

# 6) "Neither Agree Nor Disagree" (removes gramStart/gramEnd proofErr around "Nor")
$d.Content.Find.Execute("Neither Agree Nor Disagree", $false, $false, $false, $false, $false, $true, 1, $false, "Neither Agree Nor Disagree", 2)

# 7) "Here Samuel Yambo is reviewing..." (removes proofErr around "Yambo")
$d.Content.Find.Execute("Here Samuel Yambo is reviewing the work done by Samuel Adkins, to ensure correctness and style. ", $false, $false, $false, $false, $false, $true, 1, $false, "Here Samuel Yambo is reviewing the work done by Samuel Adkins, to ensure correctness and style. ", 2)

# 8) "There are no particular OS that being used." (removes gramStart/gramEnd proofErr around "particular OS")
$d.Content.Find.Execute("There are no particular OS that being used.", $false, $false, $false, $false, $false, $true, 1, $false, "There are no particular OS that being used.", 2)
